$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Summary": update running totals after trade #12 closed
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.02   # Current Capital
$summary.Range("B4").Value = 0.02      # Total P&L $
$summary.Range("B5").Value = 0.03      # Total P&L %
$summary.Range("B6").Value = 12        # Total Trades
$summary.Range("B8").Value = 5         # Losing Trades
$summary.Range("B9").Value = 33.33     # Win Rate %

# ---------------------------------------------------------------------------
# Sheet "Strategy Status": update MarketMaking strategy row (row 4)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.02     # Capital
$status.Range("D4").Value = 12         # Trades
$status.Range("E4").Value = 0.02       # P&L $
$status.Range("F4").Value = 0.02       # P&L %
$status.Range("G4").Value = 33.33      # Win Rate %

# ---------------------------------------------------------------------------
# Append trade #12 to "All Trades" and "MarketMaking" sheets (new row 13)
# ---------------------------------------------------------------------------
$tradeSheets = @("All Trades", "MarketMaking")

foreach ($sheetName in $tradeSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A13").Value = 12

    # The Date column would otherwise be auto-parsed into a date serial;
    # a leading apostrophe keeps it as plain text, matching the source data.
    $ws.Range("B13").Value = "'2026-02-17"
    $ws.Range("C13").Value = "04:07:11"

    $ws.Range("D13").Value = "MarketMaking"
    $ws.Range("E13").Value = "UP"
    $ws.Range("F13").Value = 0.21
    $ws.Range("G13").Value = 0.2
    $ws.Range("H13").Value = "CLOSED"
    $ws.Range("I13").Value = -4.7619
    $ws.Range("J13").Value = -0.01
    $ws.Range("K13").Value = 100.02
    $ws.Range("L13").Value = 0
    $ws.Range("M13").Value = 0
    $ws.Range("N13").Value = 0.6
    $ws.Range("O13").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P13").Value = "early_exit"
    $ws.Range("Q13").Value = 0.11
}
